# Updates species-observation rows: rounds the Ost/Nord (Q/R) coordinates to
# whole metres on every data row, and re-applies the (reshuffled) taxon data
# for rows 3,5,6,9-17 to match the corrected source export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 494354
$ws.Range("R2").Value = 6928891

$ws.Range("A3").Value = 111868975
$ws.Range("B3").Value = 90678
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 4366
$ws.Range("F3").Value = 'Skarp dropptaggsvamp'
$ws.Range("G3").Value = 'Hydnellum peckii'
$ws.Range("H3").Value = 'Banker'
$ws.Range("Q3").Value = 494341
$ws.Range("R3").Value = 6928940
$ws.Range("S3").Value = 30

$ws.Range("Q4").Value = 494301
$ws.Range("R4").Value = 6928922

$ws.Range("A5").Value = 111870057
$ws.Range("B5").Value = 90710
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 5449
$ws.Range("F5").Value = 'Svart taggsvamp'
$ws.Range("G5").Value = 'Phellodon niger'
$ws.Range("H5").Value = '(Fr.:Fr.) P.Karst.'
$ws.Range("P5").Value = 'Motjärnen (Motjärnen), Jmt'
$ws.Range("Q5").Value = 494314
$ws.Range("R5").Value = 6928937
$ws.Range("S5").Value = 20

$ws.Range("A6").Value = 111868497
$ws.Range("B6").Value = 90666
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 4364
$ws.Range("F6").Value = 'Dropptaggsvamp'
$ws.Range("G6").Value = 'Hydnellum ferrugineum'
$ws.Range("H6").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q6").Value = 494354
$ws.Range("R6").Value = 6928891
$ws.Range("Z6").Value = '14:23'
$ws.Range("AB6").Value = '14:23'

$ws.Range("Q7").Value = 494308
$ws.Range("R7").Value = 6928910

$ws.Range("Q8").Value = 494301
$ws.Range("R8").Value = 6928922

$ws.Range("A9").Value = 111870906
$ws.Range("B9").Value = 77515
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = 'Garnlav'
$ws.Range("G9").Value = 'Alectoria sarmentosa'
$ws.Range("H9").Value = '(Ach.) Ach.'
$ws.Range("Q9").Value = 494330
$ws.Range("R9").Value = 6928848

$ws.Range("A10").Value = 111869523
$ws.Range("B10").Value = 56543
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 103021
$ws.Range("F10").Value = 'Talltita'
$ws.Range("G10").Value = 'Poecile montanus'
$ws.Range("H10").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("Q10").Value = 494333
$ws.Range("R10").Value = 6928943
$ws.Range("Z10").Value = '15:06'
$ws.Range("AB10").Value = '15:06'

$ws.Range("A11").Value = 111868438
$ws.Range("B11").Value = 89369
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 5447
$ws.Range("F11").Value = 'Vedticka'
$ws.Range("G11").Value = 'Fuscoporia viticola'
$ws.Range("H11").Value = '(Schwein.) Murrill'
$ws.Range("Q11").Value = 494363
$ws.Range("R11").Value = 6928873
$ws.Range("S11").Value = 30

$ws.Range("A12").Value = 111870880
$ws.Range("B12").Value = 90682
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 2059
$ws.Range("F12").Value = 'Skrovlig taggsvamp'
$ws.Range("G12").Value = 'Hydnellum scabrosum'
$ws.Range("H12").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("Q12").Value = 494330
$ws.Range("R12").Value = 6928848
$ws.Range("S12").Value = 20

$ws.Range("A13").Value = 111870913
$ws.Range("B13").Value = 90666
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 4364
$ws.Range("F13").Value = 'Dropptaggsvamp'
$ws.Range("G13").Value = 'Hydnellum ferrugineum'
$ws.Range("H13").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q13").Value = 494330
$ws.Range("R13").Value = 6928848
$ws.Range("S13").Value = 20

$ws.Range("A14").Value = 111869281
$ws.Range("B14").Value = 90670
$ws.Range("D14").Value = 'VU'
$ws.Range("E14").Value = 4365
$ws.Range("F14").Value = 'Smalfotad taggsvamp'
$ws.Range("G14").Value = 'Hydnellum gracilipes'
$ws.Range("H14").Value = '(P.Karst) P.Karst'
$ws.Range("Q14").Value = 494333
$ws.Range("R14").Value = 6928943
$ws.Range("S14").Value = 30
$ws.Range("Z14").Value = '15:06'
$ws.Range("AB14").Value = '15:06'

$ws.Range("A15").Value = 111868823
$ws.Range("B15").Value = 90666
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 4364
$ws.Range("F15").Value = 'Dropptaggsvamp'
$ws.Range("G15").Value = 'Hydnellum ferrugineum'
$ws.Range("H15").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("P15").Value = 'Kläppberget, Kläppberget, Haverö, Jmt'
$ws.Range("Q15").Value = 494338
$ws.Range("R15").Value = 6928937
$ws.Range("S15").Value = 25
$ws.Range("Z15").Value = '14:23'
$ws.Range("AB15").Value = '14:23'

$ws.Range("A16").Value = 111868443
$ws.Range("B16").Value = 94134
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 53
$ws.Range("F16").Value = 'Vedtrappmossa'
$ws.Range("G16").Value = 'Crossocalyx hellerianus'
$ws.Range("H16").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("Q16").Value = 494363
$ws.Range("R16").Value = 6928873
$ws.Range("S16").Value = 30

$ws.Range("A17").Value = 111870556
$ws.Range("B17").Value = 89768
$ws.Range("D17").Value = 'VU'
$ws.Range("E17").Value = 298
$ws.Range("F17").Value = 'Laxgröppa'
$ws.Range("G17").Value = 'Byssomerulius albostramineus'
$ws.Range("H17").Value = '(Torrend) Hjortstam'
$ws.Range("Q17").Value = 494301
$ws.Range("R17").Value = 6928922
$ws.Range("S17").Value = 20
